$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7476077079772949
$ws.Range("B1").Value = 1.576925992965698
$ws.Range("C1").Value = 4.693437576293945
$ws.Range("D1").Value = 2.396350145339966
$ws.Range("E1").Value = 1.265293121337891
